# Auto update Excel log
# Appends new sensor-log rows to the "Proximity" sheet (rows 61-67) and the
# "Camera" sheet (rows 39-42), matching the data captured on 2026-02-01.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Proximity sheet: Date | Timestamp | Hour | Location | Value | Status(msg)
# ---------------------------------------------------------------------------
$proximity = $wb.Worksheets.Item("Proximity")

# Pre-format column A for the new rows as Text so Excel does not silently
# convert the "YYYY-MM-DD" strings into date serial numbers.
$proximity.Range("A61:A67").NumberFormat = "@"

$proximityRows = @(
    @("2026-02-01", "14:28:42", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "14:28:45", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "14:28:48", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "14:28:51", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "14:29:11", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "14:29:20", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "14:29:40", "14:00", "Bathroom Door",         "ENTER", "User ENTERED Bathroom")
)

$startRow = 61
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $r = $startRow + $i
    $row = $proximityRows[$i]
    $proximity.Cells.Item($r, 1).Value = $row[0]
    $proximity.Cells.Item($r, 2).Value = $row[1]
    $proximity.Cells.Item($r, 3).Value = $row[2]
    $proximity.Cells.Item($r, 4).Value = $row[3]
    $proximity.Cells.Item($r, 5).Value = $row[4]
    $proximity.Cells.Item($r, 6).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Camera sheet: Date | Timestamp | Hour | Location | Value | Status
# ---------------------------------------------------------------------------
$camera = $wb.Worksheets.Item("Camera")

$camera.Range("A39:A42").NumberFormat = "@"

$cameraRows = @(
    @("2026-02-01", "14:28:44", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "14:28:50", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "14:29:20", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "14:29:27", "14:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow = 39
for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $r = $startRow + $i
    $row = $cameraRows[$i]
    $camera.Cells.Item($r, 1).Value = $row[0]
    $camera.Cells.Item($r, 2).Value = $row[1]
    $camera.Cells.Item($r, 3).Value = $row[2]
    $camera.Cells.Item($r, 4).Value = $row[3]
    $camera.Cells.Item($r, 5).Value = $row[4]
    $camera.Cells.Item($r, 6).Value = $row[5]
}
